$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 270, shifting rows 270:326+ down by one.
$ws.Rows.Item(270).Insert()

# The inserted row inherits the formatting of the row above it (row 269,
# which is wrapped/centered text format) for every column. Column B here
# needs to hold a genuine number (General format) and column C needs the
# left-aligned wrap-text "note" format, so pull those two cell formats in
# from existing rows that already use them before writing the values.
$ws.Range("B248").Copy()
$ws.Range("B270").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C190").Copy()
$ws.Range("C270").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Application.CutCopyMode = $false

# Populate the newly inserted row 270 with the new record.
$ws.Cells.Item(270, 1).Value = "奇美實業股份有限公司"
$ws.Cells.Item(270, 2).Value = 68387705
$ws.Cells.Item(270, 3).Value = "PER MT 單位 TNE"

# Match the author's final view/selection state.
$ws.Application.Goto($ws.Range("A270"), $true)
$ws.Range("A270").Select()
